# Rename the inline picture shapes (logo images) in the document's
# footers/header, as per the authoring change:
#   footer (default)    : image1.png -> image2.png   (Pearson logo, docPr id=1)
#   footer (first page) : image1.png -> image2.png   (Pearson logo, docPr id=2)
#   header (first page) : image2.jpg -> image1.jpg   (BTec logo,    docPr id=3)

$d = $word.ActiveDocument
$sec = $d.Sections(1)

# --- Default (primary) footer: Pearson logo, id=1, image1.png -> image2.png
$footerDefault = $sec.Footers(1)
if ($footerDefault.Exists) {
    for ($i = 1; $i -le $footerDefault.Range.InlineShapes.Count; $i++) {
        $shp = $footerDefault.Range.InlineShapes($i)
        if ($shp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
            $shp.Name = "image2.png"
        }
    }
}
Write-Host "Updated default footer logo"

# --- First-page footer: Pearson logo, id=2, image1.png -> image2.png
$footerFirst = $sec.Footers(2)
if ($footerFirst.Exists) {
    for ($i = 1; $i -le $footerFirst.Range.InlineShapes.Count; $i++) {
        $shp = $footerFirst.Range.InlineShapes($i)
        if ($shp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
            $shp.Name = "image2.png"
        }
    }
}
Write-Host "Updated first-page footer logo"

# --- First-page header: BTec logo, id=3, image2.jpg -> image1.jpg
$headerFirst = $sec.Headers(2)
if ($headerFirst.Exists) {
    for ($i = 1; $i -le $headerFirst.Range.InlineShapes.Count; $i++) {
        $shp = $headerFirst.Range.InlineShapes($i)
        if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
            $shp.Name = "image1.jpg"
        }
    }
}
Write-Host "Updated first-page header logo"
